# "added length to the line"
#
# The sheet tracks individual work sessions (year/month/day/start/end/
# minutes/hours) in rows 2-148, followed by a blank spacer row and three
# summary rows (sum in minutes, sum in hours, sum in working weeks).
#
# This change adds one more work session - 2014-07-23, 19:00-22:00 (3h) -
# as a new row 149, which pushes the spacer + summary rows down by one
# row (now rows 150-153) and extends the running-total formula to cover
# the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 149; this shifts the old row 149 (blank spacer)
# and the old summary rows 150-152 down to 150-153, inheriting their
# formatting from the row above (same as Excel's default insert behavior).
$ws.Rows("149:149").Insert()

# Populate the new data row: 2014-07-23, start 19:00, end 22:00.
$ws.Range("A149").Value = 2014
$ws.Range("B149").Value = 7
$ws.Range("C149").Value = 23
$ws.Range("D149").Value = 0.79166666666666663
$ws.Range("E149").Value = 0.91666666666666663
$ws.Range("F149").Formula = "=(E149-D149)*24*60"
$ws.Range("G149").Formula = "=F149/60"

# The "sum [min]" row (now at 151) needs its SUM range extended to
# include the new row 149.
$ws.Range("F151").Formula = "=SUM(F2:F149)"

# Match the committed selection state (row 149, column F).
$ws.Range("F149").Select()
